$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header label ---
$ws.Range("L1").Value = "MasterSheet RowNo."

# --- Fill in TotalConfirmedNewCases (G) and TotalNewDeaths (I) for existing rows 2-10 ---
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 9
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 11

$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("I10").Value = 0

# --- Add new row 11 ---
$ws.Range("A11").Value = 71
$ws.Range("B11").Value = 125
$ws.Range("C11").Value = "SUB-SAHARAN AFRICA                 "
$ws.Range("D11").Value = 43921
$ws.Range("E11").Value = "Madagascar"
$ws.Range("F11").Value = 46
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "Imported cases only"
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 5397

# --- Column widths for A:O ---
$ws.Range("A1:O1").ColumnWidth = 27

# --- Number format for date column ---
$dateCol = $ws.Range("D1:D11")
$dateCol.NumberFormat = "yyyy-mm-dd;"

# --- Alignment: center/center across the whole used range A1:O11 ---
$all = $ws.Range("A1:O11")
$all.HorizontalAlignment = -4108
$all.VerticalAlignment = -4108
